# Generate Report for Handoff
# Refresh the "Handback"/"Handoff" timestamps recorded in the localization
# status report: rows whose latest handback/handoff time was 12:21:xx or
# 12:22:1x are updated to the newly observed handback time(s).

$wb = $excel.ActiveWorkbook

# --- Overview sheet (column D = "Latest Handoff Date") ---
$ws = $wb.Sheets.Item("Overview")
$ws.Range("D7").Value  = "2016-22-19 12:22:34"
$ws.Range("D10").Value = "2016-22-19 12:22:34"
$ws.Range("D11").Value = "2016-22-19 12:22:34"
$ws.Range("D12").Value = "2016-22-19 12:22:34"
$ws.Range("D13").Value = "2016-22-19 12:22:34"
$ws.Range("D14").Value = "2016-22-19 12:22:34"
$ws.Range("D15").Value = "2016-22-19 12:22:34"
$ws.Range("D16").Value = "2016-22-19 12:22:34"

# --- zh-cn sheet (column E = "Latest Handoff Datetime") ---
$ws = $wb.Sheets.Item("zh-cn")
$ws.Range("E7").Value  = "2016-03-19 12:22:31"
$ws.Range("E10").Value = "2016-03-19 12:22:31"
$ws.Range("E11").Value = "2016-03-19 12:22:31"
$ws.Range("E12").Value = "2016-03-19 12:22:31"
$ws.Range("E13").Value = "2016-03-19 12:22:31"
$ws.Range("E14").Value = "2016-03-19 12:22:31"
$ws.Range("E15").Value = "2016-03-19 12:22:31"
$ws.Range("E16").Value = "2016-03-19 12:22:31"

# --- de-de sheet (column E = "Latest Handoff Datetime") ---
$ws = $wb.Sheets.Item("de-de")
$ws.Range("E7").Value  = "2016-03-19 12:22:34"
$ws.Range("E10").Value = "2016-03-19 12:22:34"
$ws.Range("E11").Value = "2016-03-19 12:22:34"
$ws.Range("E12").Value = "2016-03-19 12:22:34"
$ws.Range("E13").Value = "2016-03-19 12:22:34"
$ws.Range("E14").Value = "2016-03-19 12:22:34"
$ws.Range("E15").Value = "2016-03-19 12:22:34"
$ws.Range("E16").Value = "2016-03-19 12:22:34"
